$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 7150.3

# Clear Scope ID # value
$ws.Range("G10").Value = ""

# Update per-line-item pricing (H column)
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55
$ws.Range("H18").Value = 476.4
$ws.Range("H19").Value = 2143.8
$ws.Range("H20").Value = 476.4
$ws.Range("H21").Value = 476.4
$ws.Range("H22").Value = 476.4
$ws.Range("H23").Value = 2143.8

# Update TOTAL
$ws.Range("H24").Value = 7150.3
